$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 21.95096015930176
$ws.Range("D2").Value = 26457

$ws.Range("C3").Value = 18.21589469909668

$ws.Range("C4").Value = 17.30036735534668

$ws.Range("C5").Value = 17.9598331451416

$ws.Range("C6").Value = 17.35973358154297
